$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume (E) columns to Text format so that numeric-looking
# strings (e.g. "1.001", "0.4883") are preserved exactly as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.804.47'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '1.941.19'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '242.62'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.4883'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").Value = '19.47'
$ws.Range("E10").Value = '  +1.98%  '
$ws.Range("D11").Value = '106.11'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").Value = '1.942.73'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = '0.07724'
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("D14").Value = '5.361'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '275.80'
$ws.Range("E16").Value = '  -3.10%  '
$ws.Range("D17").Value = '30.808.17'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '0.000007719'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '2.196.50'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '5.485'
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '6.548'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("D25").Value = '9.731'
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("D26").Value = '167.11'
$ws.Range("E26").Value = '  -1.28%  '
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").Value = '2.164'
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("D29").Value = '0.1045'
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("D30").Value = '1.392'
$ws.Range("E30").Value = '  -3.25%  '
$ws.Range("D31").Value = '4.566'
$ws.Range("E31").Value = '  -4.06%  '
$ws.Range("D32").Value = '1.554'
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("D33").Value = '4.370'
$ws.Range("D34").Value = '0.04855'
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("D35").Value = '0.7539'
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").Value = '1.159'
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D38").Value = '2.734'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").Value = '2.659'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("D41").Value = '6.536'
$ws.Range("E41").Value = '  +1.31%  '
$ws.Range("D42").Value = '77.83'
$ws.Range("E42").Value = '  +6.96%  '
$ws.Range("D43").Value = '2.100'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '0.9081'
$ws.Range("D45").Value = '108.27'
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").Value = '0.4405'
$ws.Range("E46").Value = '  -1.51%  '
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '7.754'
$ws.Range("E48").Value = '  +3.14%  '
$ws.Range("D49").Value = '1.002.09'
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").Value = '9.295'
$ws.Range("E51").Value = '  -0.68%  '
